$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 777.1039
$ws.Range("I15").Value = 777.1039
$ws.Range("K15").Value = 2331.3117
$ws.Range("M15").Value = -2162.3117
$ws.Range("H19").Value = 6406.5835
$ws.Range("I19").Value = 4994.875
$ws.Range("J19").Value = 9230
$ws.Range("K19").Value = 4994.875
$ws.Range("L19").Value = 9230
$ws.Range("M19").Value = -4819.875
$ws.Range("N19").Value = -9580
$ws.Range("H28").Value = 1833.2727
$ws.Range("I28").Value = 745.75
$ws.Range("K28").Value = 745.75
$ws.Range("M28").Value = -260.75
$ws.Range("H29").Value = 4299.5
$ws.Range("I29").Value = 4750
$ws.Range("J29").Value = 4149.3335
$ws.Range("K29").Value = 14250
$ws.Range("L29").Value = 12448.0005
$ws.Range("M29").Value = -13969
$ws.Range("N29").Value = -13010.0005
$ws.Range("H38").Value = 4515.75
$ws.Range("I38").Value = 3770.5
$ws.Range("J38").Value = 6006.25
$ws.Range("K38").Value = 11311.5
$ws.Range("L38").Value = 18018.75
$ws.Range("M38").Value = -10939.5
$ws.Range("N38").Value = -18762.75
$ws.Range("H43").Value = 1003300.4
$ws.Range("J43").Value = 2503001.5
$ws.Range("L43").Value = 2503001.5
$ws.Range("N43").Value = -2503139.5
$ws.Range("H86").Value = 6311.1113
$ws.Range("I86").Value = 7400
$ws.Range("K86").Value = 7400
$ws.Range("M86").Value = -6277
$ws.Range("H89").Value = 6311.1113
$ws.Range("I89").Value = 7400
$ws.Range("K89").Value = 37000
$ws.Range("M89").Value = -31384
$ws.Range("H93").Value = 46601
$ws.Range("J93").Value = 46601
$ws.Range("L93").Value = 46601
$ws.Range("N93").Value = -51593
$ws.Range("H106").Value = 2579.6667
$ws.Range("I106").Value = 2388
$ws.Range("K106").Value = 2388
$ws.Range("M106").Value = -1757
$ws.Range("H116").Value = 4561.722
$ws.Range("J116").Value = 4245.1113
$ws.Range("L116").Value = 4245.1113
$ws.Range("N116").Value = -11129.1113
$ws.Range("H132").Value = 3340.4827
$ws.Range("I132").Value = 3213.7856
$ws.Range("K132").Value = 9641.356800000001
$ws.Range("M132").Value = -7111.356800000001
$ws.Range("H137").Value = 3242.762
$ws.Range("I137").Value = 3531.75
$ws.Range("J137").Value = 2318
$ws.Range("K137").Value = 10595.25
$ws.Range("L137").Value = 6954
$ws.Range("M137").Value = -8045.25
$ws.Range("N137").Value = -12054
$ws.Range("H138").Value = 8134784
$ws.Range("I138").Value = 1418
$ws.Range("J138").Value = 12827111
$ws.Range("K138").Value = 4254
$ws.Range("L138").Value = 38481333
$ws.Range("M138").Value = 886
$ws.Range("N138").Value = -38491613

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11117219
$ws.Range("I32").Value = 14495122
$ws.Range("J32").Value = 18393.953
$ws.Range("K32").Value = 14495122
$ws.Range("L32").Value = 18393.953
$ws.Range("M32").Value = -14494835
$ws.Range("N32").Value = -18967.953
$ws.Range("H63").Value = 5075.231
$ws.Range("I63").Value = 3261.1052
$ws.Range("K63").Value = 3261.1052
$ws.Range("M63").Value = -2575.1052
$ws.Range("H66").Value = 5075.231
$ws.Range("I66").Value = 3261.1052
$ws.Range("K66").Value = 16305.526
$ws.Range("M66").Value = -12873.526
$ws.Range("H122").Value = 4505.1377
$ws.Range("I122").Value = 4024.9167
$ws.Range("K122").Value = 12074.7501
$ws.Range("M122").Value = -9624.750100000001
$ws.Range("H134").Value = 24999
$ws.Range("J134").Value = 24999
$ws.Range("L134").Value = 24999
$ws.Range("N134").Value = -35139

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 115
$ws.Range("I22").Value = 115
$ws.Range("K22").Value = 115
$ws.Range("M22").Value = 58
$ws.Range("H80").Value = 2208.6667
$ws.Range("J80").Value = 3770.4443
$ws.Range("L80").Value = 3770.4443
$ws.Range("N80").Value = -5766.4443
$ws.Range("H83").Value = 2208.6667
$ws.Range("J83").Value = 3770.4443
$ws.Range("L83").Value = 18852.2215
$ws.Range("N83").Value = -28836.2215
$ws.Range("H134").Value = 2684.75
$ws.Range("I134").Value = 2582.9395
$ws.Range("K134").Value = 7748.818499999999
$ws.Range("M134").Value = -5213.818499999999

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 471
$ws.Range("I7").Value = 364.33334
$ws.Range("J7").Value = 524.3333
$ws.Range("K7").Value = 364.33334
$ws.Range("L7").Value = 524.3333
$ws.Range("M7").Value = -251.33334
$ws.Range("N7").Value = -750.3333
$ws.Range("H31").Value = 21744192
$ws.Range("I31").Value = 4563.5713
$ws.Range("K31").Value = 4563.5713
$ws.Range("M31").Value = -4268.5713
$ws.Range("H34").Value = 21744192
$ws.Range("I34").Value = 4563.5713
$ws.Range("K34").Value = 4563.5713
$ws.Range("M34").Value = -4361.5713
$ws.Range("H86").Value = 4470.25
$ws.Range("I86").Value = 4293.6665
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 4293.6665
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -3170.6665
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 4470.25
$ws.Range("I89").Value = 4293.6665
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 21468.3325
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -15852.3325
$ws.Range("N89").Value = -36232
$ws.Range("H95").Value = 30623.5
$ws.Range("J95").Value = 30623.5
$ws.Range("L95").Value = 30623.5
$ws.Range("N95").Value = -36115.5
$ws.Range("H107").Value = 1246.125
$ws.Range("I107").Value = 1248.4445
$ws.Range("K107").Value = 1248.4445
$ws.Range("M107").Value = 671.5554999999999
$ws.Range("H134").Value = 1246.3043
$ws.Range("I134").Value = 1246.3043
$ws.Range("K134").Value = 3738.9129
$ws.Range("M134").Value = -1203.9129

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2062.111
$ws.Range("I5").Value = 937.1429000000001
$ws.Range("K5").Value = 2811.4287
$ws.Range("M5").Value = -2699.4287
$ws.Range("H34").Value = 1826.72
$ws.Range("J34").Value = 2615.2354
$ws.Range("L34").Value = 7845.706200000001
$ws.Range("N34").Value = -8013.706200000001
$ws.Range("H55").Value = 6411056.5
$ws.Range("J55").Value = 3175525.5
$ws.Range("L55").Value = 9526576.5
$ws.Range("N55").Value = -9526930.5
$ws.Range("H132").Value = 1451929.5
$ws.Range("I132").Value = 2286.3635
$ws.Range("K132").Value = 20577.2715
$ws.Range("M132").Value = -18047.2715
$ws.Range("H135").Value = 2062.111
$ws.Range("I135").Value = 937.1429000000001
$ws.Range("K135").Value = 8434.286100000001
$ws.Range("M135").Value = -5899.286100000001

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 44950
$ws.Range("J63").Value = 59900
$ws.Range("L63").Value = 59900
$ws.Range("N63").Value = -61272
$ws.Range("H66").Value = 44950
$ws.Range("J66").Value = 59900
$ws.Range("L66").Value = 179700
$ws.Range("N66").Value = -186564
$ws.Range("H80").Value = 2995
$ws.Range("I80").Value = 2995
$ws.Range("K80").Value = 2995
$ws.Range("M80").Value = -1997
$ws.Range("H83").Value = 2995
$ws.Range("I83").Value = 2995
$ws.Range("K83").Value = 14975
$ws.Range("M83").Value = -9983
$ws.Range("H92").Value = 25678.572
$ws.Range("J92").Value = 25678.572
$ws.Range("L92").Value = 25678.572
$ws.Range("N92").Value = -29422.572
$ws.Range("H122").Value = 1265
$ws.Range("I122").Value = 397.5
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 1192.5
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = 1257.5
$ws.Range("N122").Value = -13900

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1130.3334
$ws.Range("I16").Value = 1130.3334
$ws.Range("K16").Value = 1130.3334
$ws.Range("M16").Value = -960.3334
$ws.Range("H22").Value = 2362.8948
$ws.Range("I22").Value = 1742.5714
$ws.Range("J22").Value = 2724.75
$ws.Range("K22").Value = 1742.5714
$ws.Range("L22").Value = 2724.75
$ws.Range("M22").Value = -1447.5714
$ws.Range("N22").Value = -3314.75
$ws.Range("H27").Value = 2362.8948
$ws.Range("I27").Value = 1742.5714
$ws.Range("J27").Value = 2724.75
$ws.Range("K27").Value = 1742.5714
$ws.Range("L27").Value = 2724.75
$ws.Range("M27").Value = -1635.5714
$ws.Range("N27").Value = -2938.75
$ws.Range("H46").Value = 1235.52
$ws.Range("I46").Value = 671.9
$ws.Range("K46").Value = 671.9
$ws.Range("M46").Value = -483.9
$ws.Range("H55").Value = 632.8095
$ws.Range("J55").Value = 999.8889
$ws.Range("L55").Value = 999.8889
$ws.Range("N55").Value = -1345.8889
$ws.Range("H93").Value = 2112.4443
$ws.Range("J93").Value = 3000
$ws.Range("L93").Value = 3000
$ws.Range("N93").Value = -5496
$ws.Range("H132").Value = 254547090
$ws.Range("I132").Value = 1134.8572
$ws.Range("K132").Value = 3404.5716
$ws.Range("M132").Value = -874.5715999999998

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 126252240
$ws.Range("I100").Value = 144288060
$ws.Range("K100").Value = 288576120
$ws.Range("M100").Value = -288575579
$ws.Range("H122").Value = 41668224
$ws.Range("I122").Value = 43479690
$ws.Range("K122").Value = 130439070
$ws.Range("M122").Value = -130436620
